$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new test case row 31 (DRAIAM111 / OP11 / LogIn / Y) --------------
# Start from row 30's formatting (border + fill/wrap pattern) so the new
# row matches the look of the existing data rows, then overwrite the
# values with the new test case's data.
$ws.Range("A30:E30").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A31").Value = "DRAIAM111"
$ws.Range("B31").Value = "OP11"
$ws.Range("C31").Value = "LogIn"
$ws.Range("D31").Value = "Y"

# --- Move the view/selection down onto the newly-added row ----------------
$ws.Range("C31").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 3
